$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append new log row 60 -------------------------------------------------
$ws.Range("A60").Value = "Retour aanmelden"
$ws.Range("B60").Value = "mailmind.test@zohomail.eu"
$ws.Range("C60").Value = "Ik wil graag een artikel retourneren. Hoe werkt dat?"
$ws.Range("D60").Value = "Retour / Terugbetaling"
$ws.Range("E60").Value = "Beste klant,`nBedankt voor je bericht. Om een artikel te retourneren, verzoeken we je om contact met ons op te nemen via onze klantenservice of het retourformulier op onze website in te vullen. Vermeld hierbij het ordernummer en de reden van retour. `nZodra we deze informatie hebben ontvangen, nemen we contact met je op om de retourzending verder af te handelen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$ws.Range("F60").Value = "2025-06-22 22:17:48"
$ws.Range("G60").Value = "Ja"

# Reset the row height back to the sheet default — entering multi-line text
# via .Value auto-wraps/sizes the row; AutoFit drops the custom-height flag
# so row 60 serializes the same way the other rows do.
$ws.Rows.Item(60).AutoFit()

# --- Extend conditional formatting ranges to include the new row ----------
$fcsD = $ws.Range("D2:D59").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($ws.Range("D2:D60"))
}

$fcsG = $ws.Range("G2:G59").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($ws.Range("G2:G60"))
}

# --- Update Dashboard summary count for "Retour / Terugbetaling" ----------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 9
